$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(23).Clear()
$ws.Rows(24).Clear()
$ws.Columns("BA").Clear()
$ws.Range("B1").Value = 39583
$ws.Range("C1").Value = 39765
$ws.Range("D1").Value = 39948
$ws.Range("E1").Value = 40130
$ws.Range("F1").Value = 40310
$ws.Range("G1").Value = 40494
$ws.Range("H1").Value = 40676
$ws.Range("I1").Value = 40862
$ws.Range("J1").Value = 41044
$ws.Range("K1").Value = 41228
$ws.Range("L1").Value = 41409
$ws.Range("M1").Value = 41592
$ws.Range("N1").Value = 41774
$ws.Range("O1").Value = 41957
$ws.Range("P1").Value = 42137
$ws.Range("Q1").Value = 42321
$ws.Range("R1").Value = 42503
$ws.Range("S1").Value = 42689
$ws.Range("T1").Value = 42867
$ws.Range("U1").Value = 43053
$ws.Range("V1").Value = 43145
$ws.Range("W1").Value = 43235
$ws.Range("X1").Value = 43326
$ws.Range("Y1").Value = 43418
$ws.Range("Z1").Value = 43510
$ws.Range("AA1").Value = 43600
$ws.Range("AB1").Value = 43691
$ws.Range("AC1").Value = 43783
$ws.Range("AD1").Value = 43875
$ws.Range("AE1").Value = 43966
$ws.Range("AF1").Value = 44068
$ws.Range("AG1").Value = 44159
$ws.Range("AH1").Value = 44251
$ws.Range("AI1").Value = 44341
$ws.Range("AJ1").Value = 44432
$ws.Range("AK1").Value = 44525
$ws.Range("AL1").Value = 44617
$ws.Range("AM1").Value = 44706
$ws.Range("AN1").Value = 44798
$ws.Range("AO1").Value = 44890
$ws.Range("AP1").Value = 44981
$ws.Range("AQ1").Value = 45071
$ws.Range("AR1").Value = 45163
$ws.Range("AS1").Value = 45254
$ws.Range("AT1").Value = 45345
$ws.Range("AU1").Value = 45436
$ws.Range("AV1").Value = 45534
$ws.Range("AW1").Value = 45618
$ws.Range("AX1").Value = 45713
$ws.Range("AY1").Value = 45800
$ws.Range("AZ1").Value = 45891
$ws.Range("B3").Value = 1.328558632615739
$ws.Range("C3").Value = 0.4282194198276246
$ws.Range("D3").Value = -1.435981453719049
$ws.Range("B4").Value = 1.356849765318358
$ws.Range("C4").Value = 0.5762410323606026
$ws.Range("D4").Value = -0.7704417043119083
$ws.Range("E4").Value = -0.8235211753995442
$ws.Range("F4").Value = 0.406633294022174
$ws.Range("C5").Clear()
$ws.Range("D5").Value = -0.6484173407089511
$ws.Range("E5").Value = -0.735462857474789
$ws.Range("F5").Value = -0.1034614224434405
$ws.Range("G5").Value = 1.121293995080253
$ws.Range("H5").Value = 2.185496833134781
$ws.Range("E6").Clear()
$ws.Range("F6").Value = -0.2218152533720597
$ws.Range("G6").Value = 0.1480821602630744
$ws.Range("H6").Value = 0.6652762968575532
$ws.Range("I6").Value = 1.665250327443002
$ws.Range("J6").Value = 0.8574941660507873
$ws.Range("G7").Clear()
$ws.Range("H7").Value = 0.6294804885647043
$ws.Range("I7").Value = 1.617722252039
$ws.Range("J7").Value = 1.693469135756587
$ws.Range("K7").Value = 1.079796209653616
$ws.Range("L7").Value = 0.1494732105682406
$ws.Range("I8").Clear()
$ws.Range("J8").Value = 1.77975641852226
$ws.Range("K8").Value = 1.488234279941625
$ws.Range("L8").Value = 0.8024032015999882
$ws.Range("M8").Value = 1.374377011838535
$ws.Range("N8").Value = 1.656936590801972
$ws.Range("K9").Clear()
$ws.Range("L9").Value = 0.8024032016000104
$ws.Range("M9").Value = 0.8909614188480353
$ws.Range("N9").Value = 0.922773818606859
$ws.Range("O9").Value = 1.310895847186577
$ws.Range("P9").Value = 1.346932828201242
$ws.Range("M10").Clear()
$ws.Range("N10").Value = 0.8507004532711138
$ws.Range("O10").Value = 1.17312580610518
$ws.Range("P10").Value = 1.364302026343633
$ws.Range("Q10").Value = 1.862478303083726
$ws.Range("R10").Value = 1.745747589686109
$ws.Range("O11").Clear()
$ws.Range("P11").Value = 1.342708276326299
$ws.Range("Q11").Value = 1.804078246438934
$ws.Range("R11").Value = 1.644798626926303
$ws.Range("S11").Value = 1.639776099317536
$ws.Range("T11").Value = 1.843649045891893
$ws.Range("Q12").Clear()
$ws.Range("R12").Value = 1.609625625599986
$ws.Range("S12").Value = 1.601529483008668
$ws.Range("T12").Value = 1.741128155516525
$ws.Range("U12").Value = 2.181728312936415
$ws.Range("V12").Value = 2.284406789710336
$ws.Range("W12").Value = 2.463589365374652
$ws.Range("X12").Value = 2.349806433215029
$ws.Range("R13:S13").Clear()
$ws.Range("T13").Value = 1.723022434657207
$ws.Range("U13").Value = 1.99880983009828
$ws.Range("V13").Value = 1.990690441067144
$ws.Range("W13").Value = 2.149194501693219
$ws.Range("X13").Value = 2.036910005299108
$ws.Range("Y13").Value = 2.010025322622599
$ws.Range("Z13").Value = 1.665971362160357
$ws.Range("AA13").Value = 1.332860091726285
$ws.Range("AB13").Value = 1.029194292875912
$ws.Range("T14:V14").Clear()
$ws.Range("W14").Value = 2.104406008906734
$ws.Range("X14").Value = 2.104406008906734
$ws.Range("Y14").Value = 2.088155108730527
$ws.Range("Z14").Value = 2.031292234149706
$ws.Range("AA14").Value = 1.799885362733189
$ws.Range("AB14").Value = 1.31420459445093
$ws.Range("AC14").Value = 0.7771393814490102
$ws.Range("AD14").Value = 0.4126128934655471
$ws.Range("AE14").Value = 0.2336391425753925
$ws.Range("AF14").Value = -4.43626840667447
$ws.Range("V15:Z15").Clear()
$ws.Range("AA15").Value = 1.893295577996756
$ws.Range("AB15").Value = 1.650113596657588
$ws.Range("AC15").Value = 1.407367165006201
$ws.Range("AD15").Value = 1.156986202028509
$ws.Range("AE15").Value = 0.9207450904090253
$ws.Range("AF15").Value = -2.63419394755392
$ws.Range("AG15").Value = -2.71887004062904
$ws.Range("AH15").Value = -2.96879819115512
$ws.Range("AI15").Value = -2.010709456685855
$ws.Range("AJ15").Value = -1.513408827666285
$ws.Range("Y16:AD16").Clear()
$ws.Range("AE16").Value = 0.9554040357173665
$ws.Range("AF16").Value = -1.6615457433243
$ws.Range("AG16").Value = -1.937757788996253
$ws.Range("AH16").Value = -2.438555173006141
$ws.Range("AI16").Value = -1.14257141002756
$ws.Range("AJ16").Value = 0.7106578563214505
$ws.Range("AK16").Value = 0.4582698374457683
$ws.Range("AL16").Value = 1.154413086110817
$ws.Range("AM16").Value = 1.5286818008164
$ws.Range("AN16").Value = 1.618732201786743
$ws.Range("AC17:AG17").Clear()
$ws.Range("AH17").Value = -2.413633125962611
$ws.Range("AI17").Value = -1.976964192800379
$ws.Range("AJ17").Value = -1.416208084524317
$ws.Range("AK17").Value = -1.590220060268321
$ws.Range("AL17").Value = -1.166698219025086
$ws.Range("AM17").Value = -0.7118141543333012
$ws.Range("AN17").Value = -0.4630595634534385
$ws.Range("AO17").Value = 1.314675624401973
$ws.Range("AP17").Value = 0.006126408955742235
$ws.Range("AQ17").Value = 0.001079933351455509
$ws.Range("AR17").Value = -0.09609276733164585
$ws.Range("AG18:AK18").Clear()
$ws.Range("AL18").Value = -1.366617752737886
$ws.Range("AM18").Value = -1.247796181364325
$ws.Range("AN18").Value = -1.271569079498103
$ws.Range("AO18").Value = 3.40905661301254
$ws.Range("AP18").Value = 0.5370151562237302
$ws.Range("AQ18").Value = 0.6889047703476203
$ws.Range("AR18").Value = 0.5285660612534882
$ws.Range("AS18").Value = 0.1199358335146838
$ws.Range("AT18").Value = -0.1634698065940632
$ws.Range("AU18").Value = -0.00209793826797533
$ws.Range("AV18").Value = -0.02761034355766023
$ws.Range("AK19:AO19").Clear()
$ws.Range("AP19").Value = 0.6286476167952193
$ws.Range("AQ19").Value = 0.8380511040219529
$ws.Range("AR19").Value = 1.051341880573431
$ws.Range("AS19").Value = 0.541919990776929
$ws.Range("AT19").Value = -0.1145111565623136
$ws.Range("AU19").Value = 0.286657616500996
$ws.Range("AV19").Value = 0.08174908622293753
$ws.Range("AW19").Value = 0.2676745853112728
$ws.Range("AX19").Value = 0.4725905789402463
$ws.Range("AY19").Value = 0.501314651583451
$ws.Range("AZ19").Value = 0.5152269879013183
$ws.Range("AO20:AS20").Clear()
$ws.Range("AT20").Value = -0.1001442164906607
$ws.Range("AU20").Value = 0.2192992773568214
$ws.Range("AV20").Value = 0.01413548204556303
$ws.Range("AW20").Value = -0.05216314813395462
$ws.Range("AX20").Value = 0.08564335928031852
$ws.Range("AY20").Value = 0.2155158706220295
$ws.Range("AZ20").Value = 0.3332251551730891
$ws.Range("AS21:AW21").Clear()
$ws.Range("AX21").Value = 0.02760436504196662
$ws.Range("AY21").Value = 0.1625720590646029
$ws.Range("AZ21").Value = 0.2843066506847514
$ws.Range("AW22:AZ22").Clear()
